# Replace "OIE" with "WOAH" throughout the FMDV story-map workbook.
# The change is purely textual: every standalone occurrence of "OIE"
# (as used for the old name of the World Organisation for Animal Health)
# becomes "WOAH", while unrelated text (e.g. the lower-case "oie.int"
# URLs) must be left untouched.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("Sheet 1")
$sheet2 = $wb.Worksheets.Item("References")

# Cells on "Sheet 1" (column E - "Content") that contain the word "OIE".
$sheet1Cells = @("E5", "E6", "E7", "E14", "E17", "E21", "E42", "E54", "E66", "E92")

foreach ($addr in $sheet1Cells) {
    $cell = $sheet1.Range($addr)
    $text = $cell.Value2
    if ($text -ne $null) {
        $cell.Value = [System.Text.RegularExpressions.Regex]::Replace($text, "\bOIE\b", "WOAH")
    }
}

# Cells on "References" (column C - reference text) that contain the word "OIE".
$sheet2Cells = @("C2", "C6", "C9", "C10", "C11")

foreach ($addr in $sheet2Cells) {
    $cell = $sheet2.Range($addr)
    $text = $cell.Value2
    if ($text -ne $null) {
        $cell.Value = [System.Text.RegularExpressions.Regex]::Replace($text, "\bOIE\b", "WOAH")
    }
}
